$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-16 05:18:58'
$ws.Range('E3').Value = '2026-02-16 05:19:00'
$ws.Range('E4').Value = '2026-02-16 05:19:02'
$ws.Range('H4').Value = '''67%'
$ws.Range('J4').Value = '1014.4 hPa'
$ws.Range('L4').Value = '29.9 km/h - 322º 4:44 TU'
$ws.Range('E5').Value = '2026-02-16 05:19:05'
$ws.Range('I5').Value = '1.8 mm'
$ws.Range('M5').Value = '-0.9 °C 4:57 TU'
$ws.Range('E6').Value = '2026-02-16 05:19:07'
$ws.Range('H6').Value = '''85%'
$ws.Range('J6').Value = '1014.5 hPa'
$ws.Range('O6').Value = '6.8 °C'
$ws.Range('E7').Value = '2026-02-16 05:19:10'
$ws.Range('J7').Value = '1014.8 hPa'
$ws.Range('M7').Value = '13.7 °C 4:55 TU'
$ws.Range('O7').Value = '13.2 °C'
$ws.Range('E8').Value = '2026-02-16 05:19:12'
$ws.Range('J8').Value = '1014.6 hPa'
$ws.Range('M8').Value = '9.9 °C 4:56 TU'
$ws.Range('E9').Value = '2026-02-16 05:19:15'
$ws.Range('H9').Value = '''96%'
$ws.Range('O9').Value = '5.1 °C'
$ws.Range('E10').Value = '2026-02-16 05:19:18'
$ws.Range('E11').Value = '2026-02-16 05:19:20'
$ws.Range('N11').Value = '-0.1 °C 4:56 TU'
$ws.Range('O11').Value = '0.8 °C'
$ws.Range('E12').Value = '2026-02-16 05:19:23'
$ws.Range('N12').Value = '4.4 °C 4:51 TU'
$ws.Range('O12').Value = '5.5 °C'
$ws.Range('E13').Value = '2026-02-16 05:19:25'
$ws.Range('H13').Value = '''88%'
$ws.Range('N13').Value = '-0.6 °C 4:43 TU'
$ws.Range('O13').Value = '1.1 °C'
$ws.Range('E14').Value = '2026-02-16 05:19:27'
$ws.Range('O14').Value = '12.8 °C'
$ws.Range('E15').Value = '2026-02-16 05:19:30'
$ws.Range('N15').Value = '3.3 °C 4:52 TU'
$ws.Range('O15').Value = '5.2 °C'
$ws.Range('E16').Value = '2026-02-16 05:19:32'
$ws.Range('I16').Value = '1.7 mm'
$ws.Range('L16').Value = '90.7 km/h - 225º 4:59 TU'
$ws.Range('M16').Value = '0.7 °C 4:59 TU'
$ws.Range('O16').Value = '-0.6 °C'
$ws.Range('E17').Value = '2026-02-16 05:19:35'
$ws.Range('O17').Value = '5.4 °C'
$ws.Range('E18').Value = '2026-02-16 05:19:37'
$ws.Range('J18').Value = '1014.9 hPa'
$ws.Range('N18').Value = '2.9 °C 4:58 TU'
$ws.Range('O18').Value = '4.2 °C'
$ws.Range('E19').Value = '2026-02-16 05:19:39'
$ws.Range('H19').Value = '''96%'
$ws.Range('N19').Value = '2.5 °C 4:57 TU'
$ws.Range('E20').Value = '2026-02-16 05:19:42'
$ws.Range('H20').Value = '''92%'
$ws.Range('N20').Value = '-2.0 °C 4:45 TU'
$ws.Range('E21').Value = '2026-02-16 05:19:44'
$ws.Range('O21').Value = '4.7 °C'
$ws.Range('E22').Value = '2026-02-16 05:19:47'
$ws.Range('I22').Value = '0.7 mm'
$ws.Range('N22').Value = '-6.5 °C 4:30 TU'
$ws.Range('E23').Value = '2026-02-16 05:19:49'
$ws.Range('H23').Value = '''85%'
$ws.Range('I23').Value = '0.7 mm'
$ws.Range('M23').Value = '-0.4 °C 4:58 TU'
$ws.Range('E24').Value = '2026-02-16 05:19:52'
$ws.Range('J24').Value = '1018.0 hPa'
$ws.Range('O24').Value = '10.5 °C'
$ws.Range('E25').Value = '2026-02-16 05:19:54'
$ws.Range('E26').Value = '2026-02-16 05:19:57'
$ws.Range('E27').Value = '2026-02-16 05:19:59'
$ws.Range('N27').Value = '-0.1 °C 4:39 TU'
$ws.Range('O27').Value = '0.8 °C'
$ws.Range('E28').Value = '2026-02-16 05:20:02'
$ws.Range('H28').Value = '''91%'
$ws.Range('N28').Value = '1.8 °C 4:59 TU'
$ws.Range('O28').Value = '3.2 °C'
$ws.Range('E29').Value = '2026-02-16 05:20:04'
$ws.Range('E30').Value = '2026-02-16 05:20:07'
$ws.Range('J30').Value = '1014.5 hPa'
$ws.Range('N30').Value = '6.3 °C 4:48 TU'
$ws.Range('E31').Value = '2026-02-16 05:20:09'
$ws.Range('H31').Value = '''56%'
$ws.Range('O31').Value = '14.0 °C'
$ws.Range('E32').Value = '2026-02-16 05:20:12'
$ws.Range('H32').Value = '''81%'
$ws.Range('E33').Value = '2026-02-16 05:20:14'
$ws.Range('H33').Value = '''73%'
$ws.Range('J33').Value = '1016.1 hPa'
$ws.Range('N33').Value = '2.0 °C 4:51 TU'
$ws.Range('O33').Value = '4.4 °C'
$ws.Range('E34').Value = '2026-02-16 05:20:17'
$ws.Range('E35').Value = '2026-02-16 05:20:19'
$ws.Range('J35').Value = '1019.1 hPa'
$ws.Range('N35').Value = '6.6 °C 4:51 TU'
$ws.Range('E36').Value = '2026-02-16 05:20:22'
$ws.Range('H36').Value = '''92%'
$ws.Range('L36').Value = '10.4 km/h - 62º 4:46 TU'
$ws.Range('E37').Value = '2026-02-16 05:20:24'
$ws.Range('N37').Value = '0.7 °C 4:41 TU'
$ws.Range('O37').Value = '1.7 °C'
$ws.Range('E38').Value = '2026-02-16 05:20:27'
$ws.Range('E39').Value = '2026-02-16 05:20:30'
$ws.Range('E40').Value = '2026-02-16 05:20:32'
$ws.Range('H40').Value = '''96%'
$ws.Range('O40').Value = '2.9 °C'
$ws.Range('E41').Value = '2026-02-16 05:20:35'
$ws.Range('E42').Value = '2026-02-16 05:20:38'
$ws.Range('N42').Value = '5.7 °C 4:34 TU'
$ws.Range('O42').Value = '6.4 °C'
$ws.Range('E43').Value = '2026-02-16 05:20:40'
$ws.Range('N43').Value = '2.2 °C 4:53 TU'
$ws.Range('O43').Value = '3.4 °C'
$ws.Range('E44').Value = '2026-02-16 05:20:43'
$ws.Range('M44').Value = '0.7 °C 4:39 TU'
$ws.Range('O44').Value = '-0.2 °C'
$ws.Range('E45').Value = '2026-02-16 05:20:46'
$ws.Range('J45').Value = '1019.6 hPa'
$ws.Range('L45').Value = '5.4 km/h - 216º 4:34 TU'
$ws.Range('E46').Value = '2026-02-16 05:20:48'
$ws.Range('H46').Value = '''62%'
$ws.Range('K46').Value = '-0.1 MJ/m2'
